$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value2 = 2585.6365   # H40: 2283.3684 -> 2585.6365
$ws.Cells.Item(40, 9).Value2 = 2380.5   # I40: 1802.4546 -> 2380.5
$ws.Cells.Item(40, 11).Value2 = 2380.5   # K40: 1802.4546 -> 2380.5
$ws.Cells.Item(40, 13).Value2 = -2205.5   # M40: -1627.4546 -> -2205.5
$ws.Cells.Item(41, 8).Value2 = 1811.9   # H41: 1987.3334 -> 1811.9
$ws.Cells.Item(41, 9).Value2 = 1531.5   # I41: 1857.8 -> 1531.5
$ws.Cells.Item(41, 10).Value2 = 2232.5   # J41: 2149.25 -> 2232.5
$ws.Cells.Item(41, 11).Value2 = 1531.5   # K41: 1857.8 -> 1531.5
$ws.Cells.Item(41, 12).Value2 = 2232.5   # L41: 2149.25 -> 2232.5
$ws.Cells.Item(41, 13).Value2 = -1091.5   # M41: -1417.8 -> -1091.5
$ws.Cells.Item(41, 14).Value2 = -3112.5   # N41: -3029.25 -> -3112.5
$ws.Cells.Item(92, 8).Value2 = 640.8   # H92: 666.3 -> 640.8
$ws.Cells.Item(92, 10).Value2 = 740.3333   # J92: 825.3333 -> 740.3333
$ws.Cells.Item(92, 12).Value2 = 740.3333   # L92: 825.3333 -> 740.3333
$ws.Cells.Item(92, 14).Value2 = -3236.3333   # N92: -3321.3333 -> -3236.3333
$ws.Cells.Item(107, 8).Value2 = 5789.56   # H107: 5789.76 -> 5789.56
$ws.Cells.Item(107, 10).Value2 = 2662.4546   # J107: 2662.9092 -> 2662.4546
$ws.Cells.Item(107, 12).Value2 = 2662.4546   # L107: 2662.9092 -> 2662.4546
$ws.Cells.Item(107, 14).Value2 = -6502.4546   # N107: -6502.9092 -> -6502.4546
$ws.Cells.Item(113, 8).Value2 = 36901   # H113: 26167.334 -> 36901
$ws.Cells.Item(113, 9).Value2 = 36901   # I113: 30400.8 -> 36901
$ws.Cells.Item(113, 10).Value2 = 0   # J113: 5000 -> 0
$ws.Cells.Item(113, 11).Value2 = 36901   # K113: 30400.8 -> 36901
$ws.Cells.Item(113, 12).Value2 = 0   # L113: 5000 -> 0
$ws.Cells.Item(113, 13).ClearContents()   # M113: clear (was -27146.8)
$ws.Cells.Item(113, 14).Value2 = -33647   # N113: -11508 -> -33647
$ws.Cells.Item(136, 8).Value2 = 96999.5   # H136: 97000 -> 96999.5
$ws.Cells.Item(136, 10).Value2 = 96999.5   # J136: 97000 -> 96999.5
$ws.Cells.Item(136, 12).Value2 = 96999.5   # L136: 97000 -> 96999.5
$ws.Cells.Item(136, 14).Value2 = -107199.5   # N136: -107200 -> -107199.5
$ws.Cells.Item(138, 8).Value2 = 4881.721   # H138: 4867.516 -> 4881.721
$ws.Cells.Item(138, 10).Value2 = 7050.0513   # J138: 6973.825 -> 7050.0513
$ws.Cells.Item(138, 12).Value2 = 21150.1539   # L138: 20921.475 -> 21150.1539
$ws.Cells.Item(138, 14).Value2 = -31430.1539   # N138: -31201.475 -> -31430.1539

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value2 = 1406364.8   # H122: 1172512.4 -> 1406364.8
$ws.Cells.Item(122, 9).Value2 = 6999.5   # I122: 5166.3335 -> 6999.5
$ws.Cells.Item(122, 10).Value2 = 1756206.1   # J122: 1561627.6 -> 1756206.1
$ws.Cells.Item(122, 11).Value2 = 20998.5   # K122: 15499.0005 -> 20998.5
$ws.Cells.Item(122, 12).Value2 = 5268618.300000001   # L122: 4684882.800000001 -> 5268618.300000001
$ws.Cells.Item(122, 13).Value2 = -18548.5   # M122: -13049.0005 -> -18548.5
$ws.Cells.Item(122, 14).Value2 = -5273518.300000001   # N122: -4689782.800000001 -> -5273518.300000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value2 = 58354.555   # H82: 48307.777 -> 58354.555
$ws.Cells.Item(82, 9).Value2 = 11949.75   # I82: 9939.6 -> 11949.75
$ws.Cells.Item(82, 10).Value2 = 95478.39999999999   # J82: 96268 -> 95478.39999999999
$ws.Cells.Item(82, 11).Value2 = 11949.75   # K82: 9939.6 -> 11949.75
$ws.Cells.Item(82, 12).Value2 = 95478.39999999999   # L82: 96268 -> 95478.39999999999
$ws.Cells.Item(82, 13).Value2 = -11566.75   # M82: -9556.6 -> -11566.75
$ws.Cells.Item(82, 14).Value2 = -96244.39999999999   # N82: -97034 -> -96244.39999999999
$ws.Cells.Item(85, 8).Value2 = 58354.555   # H85: 48307.777 -> 58354.555
$ws.Cells.Item(85, 9).Value2 = 11949.75   # I85: 9939.6 -> 11949.75
$ws.Cells.Item(85, 10).Value2 = 95478.39999999999   # J85: 96268 -> 95478.39999999999
$ws.Cells.Item(85, 11).Value2 = 11949.75   # K85: 9939.6 -> 11949.75
$ws.Cells.Item(85, 12).Value2 = 95478.39999999999   # L85: 96268 -> 95478.39999999999
$ws.Cells.Item(85, 13).Value2 = -10623.75   # M85: -8613.6 -> -10623.75
$ws.Cells.Item(85, 14).Value2 = -98130.39999999999   # N85: -98920 -> -98130.39999999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value2 = 13169.5   # H58: 9700.556 -> 13169.5
$ws.Cells.Item(58, 9).Value2 = 14203.4   # I58: 10686.571 -> 14203.4
$ws.Cells.Item(58, 10).Value2 = 8000   # J58: 6249.5 -> 8000
$ws.Cells.Item(58, 11).Value2 = 14203.4   # K58: 10686.571 -> 14203.4
$ws.Cells.Item(58, 12).Value2 = 8000   # L58: 6249.5 -> 8000
$ws.Cells.Item(58, 13).Value2 = -14000.4   # M58: -10483.571 -> -14000.4
$ws.Cells.Item(58, 14).Value2 = -8406   # N58: -6655.5 -> -8406
$ws.Cells.Item(99, 8).Value2 = 58055556   # H99: 16595873 -> 58055556
$ws.Cells.Item(99, 9).Value2 = 58055556   # I99: 16595873 -> 58055556
$ws.Cells.Item(99, 11).Value2 = 58055556   # K99: 16595873 -> 58055556
$ws.Cells.Item(99, 13).Value2 = -58054058   # M99: -16594375 -> -58054058
$ws.Cells.Item(126, 8).Value2 = 58055556   # H126: 16595873 -> 58055556
$ws.Cells.Item(126, 9).Value2 = 58055556   # I126: 16595873 -> 58055556
$ws.Cells.Item(126, 11).Value2 = 174166668   # K126: 49787619 -> 174166668
$ws.Cells.Item(126, 13).Value2 = -174164198   # M126: -49785149 -> -174164198
$ws.Cells.Item(136, 8).Value2 = 13169.5   # H136: 9700.556 -> 13169.5
$ws.Cells.Item(136, 9).Value2 = 14203.4   # I136: 10686.571 -> 14203.4
$ws.Cells.Item(136, 10).Value2 = 8000   # J136: 6249.5 -> 8000
$ws.Cells.Item(136, 11).Value2 = 42610.2   # K136: 32059.713 -> 42610.2
$ws.Cells.Item(136, 12).Value2 = 24000   # L136: 18748.5 -> 24000
$ws.Cells.Item(136, 13).Value2 = -40060.2   # M136: -29509.713 -> -40060.2
$ws.Cells.Item(136, 14).Value2 = -29100   # N136: -23848.5 -> -29100

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value2 = 435666   # H5: 527262.3 -> 435666
$ws.Cells.Item(5, 9).Value2 = 993.8570999999999   # I5: 1032.2667 -> 993.8570999999999
$ws.Cells.Item(5, 10).Value2 = 1111822.6   # J5: 2500625 -> 1111822.6
$ws.Cells.Item(5, 11).Value2 = 2981.5713   # K5: 3096.800099999999 -> 2981.5713
$ws.Cells.Item(5, 12).Value2 = 3335467.8   # L5: 7501875 -> 3335467.8
$ws.Cells.Item(5, 13).Value2 = -2869.5713   # M5: -2984.800099999999 -> -2869.5713
$ws.Cells.Item(5, 14).Value2 = -3335691.8   # N5: -7502099 -> -3335691.8
$ws.Cells.Item(63, 8).Value2 = 3000   # H63: 0 -> 3000
$ws.Cells.Item(63, 10).Value2 = 3000   # J63: 0 -> 3000
$ws.Cells.Item(63, 12).Value2 = 9000   # L63: 0 -> 9000
$ws.Cells.Item(63, 14).Value2 = -10498   # N63: None -> -10498
$ws.Cells.Item(64, 8).Value2 = 2202.8   # H64: 2298.6 -> 2202.8
$ws.Cells.Item(64, 9).Value2 = 1345.6   # I64: 1540.25 -> 1345.6
$ws.Cells.Item(64, 10).Value2 = 3060   # J64: 2804.1667 -> 3060
$ws.Cells.Item(64, 11).Value2 = 4036.8   # K64: 4620.75 -> 4036.8
$ws.Cells.Item(64, 12).Value2 = 9180   # L64: 8412.500100000001 -> 9180
$ws.Cells.Item(64, 13).Value2 = -3766.8   # M64: -4350.75 -> -3766.8
$ws.Cells.Item(64, 14).Value2 = -9720   # N64: -8952.500100000001 -> -9720
$ws.Cells.Item(66, 8).Value2 = 3000   # H66: 0 -> 3000
$ws.Cells.Item(66, 10).Value2 = 3000   # J66: 0 -> 3000
$ws.Cells.Item(66, 12).Value2 = 27000   # L66: 0 -> 27000
$ws.Cells.Item(66, 14).Value2 = -34488   # N66: None -> -34488
$ws.Cells.Item(67, 8).Value2 = 2202.8   # H67: 2298.6 -> 2202.8
$ws.Cells.Item(67, 9).Value2 = 1345.6   # I67: 1540.25 -> 1345.6
$ws.Cells.Item(67, 10).Value2 = 3060   # J67: 2804.1667 -> 3060
$ws.Cells.Item(67, 11).Value2 = 4036.8   # K67: 4620.75 -> 4036.8
$ws.Cells.Item(67, 12).Value2 = 9180   # L67: 8412.500100000001 -> 9180
$ws.Cells.Item(67, 13).Value2 = -3100.8   # M67: -3684.75 -> -3100.8
$ws.Cells.Item(67, 14).Value2 = -11052   # N67: -10284.5001 -> -11052
$ws.Cells.Item(75, 8).Value2 = 2166.3333   # H75: 1338.4286 -> 2166.3333
$ws.Cells.Item(75, 9).Value2 = 500   # I75: 490 -> 500
$ws.Cells.Item(75, 10).Value2 = 2999.5   # J75: 1974.75 -> 2999.5
$ws.Cells.Item(75, 11).Value2 = 1500   # K75: 1470 -> 1500
$ws.Cells.Item(75, 12).Value2 = 8998.5   # L75: 5924.25 -> 8998.5
$ws.Cells.Item(75, 13).Value2 = -502   # M75: -472 -> -502
$ws.Cells.Item(75, 14).Value2 = -10994.5   # N75: -7920.25 -> -10994.5
$ws.Cells.Item(78, 8).Value2 = 2166.3333   # H78: 1338.4286 -> 2166.3333
$ws.Cells.Item(78, 9).Value2 = 500   # I78: 490 -> 500
$ws.Cells.Item(78, 10).Value2 = 2999.5   # J78: 1974.75 -> 2999.5
$ws.Cells.Item(78, 11).Value2 = 4500   # K78: 4410 -> 4500
$ws.Cells.Item(78, 12).Value2 = 26995.5   # L78: 17772.75 -> 26995.5
$ws.Cells.Item(78, 13).Value2 = 492   # M78: 582 -> 492
$ws.Cells.Item(78, 14).Value2 = -36979.5   # N78: -27756.75 -> -36979.5
$ws.Cells.Item(126, 8).Value2 = 12350.667   # H126: 12522.714 -> 12350.667
$ws.Cells.Item(126, 10).Value2 = 14554.8   # J126: 14388.167 -> 14554.8
$ws.Cells.Item(126, 12).Value2 = 43664.39999999999   # L126: 43164.501 -> 43664.39999999999
$ws.Cells.Item(126, 14).Value2 = -53544.39999999999   # N126: -53044.501 -> -53544.39999999999
$ws.Cells.Item(135, 8).Value2 = 435666   # H135: 527262.3 -> 435666
$ws.Cells.Item(135, 9).Value2 = 993.8570999999999   # I135: 1032.2667 -> 993.8570999999999
$ws.Cells.Item(135, 10).Value2 = 1111822.6   # J135: 2500625 -> 1111822.6
$ws.Cells.Item(135, 11).Value2 = 8944.713899999999   # K135: 9290.400299999999 -> 8944.713899999999
$ws.Cells.Item(135, 12).Value2 = 10006403.4   # L135: 22505625 -> 10006403.4
$ws.Cells.Item(135, 13).Value2 = -6409.713899999999   # M135: -6755.400299999999 -> -6409.713899999999
$ws.Cells.Item(135, 14).Value2 = -10011473.4   # N135: -22510695 -> -10011473.4

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value2 = 158.71428   # H2: 180.83333 -> 158.71428
$ws.Cells.Item(2, 9).Value2 = 101.833336   # I2: 117 -> 101.833336
$ws.Cells.Item(2, 11).Value2 = 101.833336   # K2: 117 -> 101.833336
$ws.Cells.Item(2, 13).Value2 = 11.166664   # M2: -4 -> 11.166664
$ws.Cells.Item(97, 8).Value2 = 10008.833   # H97: 9527.4 -> 10008.833
$ws.Cells.Item(97, 9).Value2 = 13019.706   # I97: 11700.263 -> 13019.706
$ws.Cells.Item(97, 10).Value2 = 2696.7144   # J97: 2646.6667 -> 2696.7144
$ws.Cells.Item(97, 11).Value2 = 13019.706   # K97: 11700.263 -> 13019.706
$ws.Cells.Item(97, 12).Value2 = 2696.7144   # L97: 2646.6667 -> 2696.7144
$ws.Cells.Item(97, 13).Value2 = -12523.706   # M97: -11204.263 -> -12523.706
$ws.Cells.Item(97, 14).Value2 = -3688.7144   # N97: -3638.6667 -> -3688.7144
$ws.Cells.Item(102, 8).Value2 = 23286   # H102: 26528.357 -> 23286
$ws.Cells.Item(102, 9).Value2 = 20851   # I102: 24534.818 -> 20851
$ws.Cells.Item(102, 10).Value2 = 33837.668   # J102: 33838 -> 33837.668
$ws.Cells.Item(102, 11).Value2 = 20851   # K102: 24534.818 -> 20851
$ws.Cells.Item(102, 12).Value2 = 33837.668   # L102: 33838 -> 33837.668
$ws.Cells.Item(102, 13).Value2 = -19229   # M102: -22912.818 -> -19229
$ws.Cells.Item(102, 14).Value2 = -37081.668   # N102: -37082 -> -37081.668
$ws.Cells.Item(122, 8).Value2 = 15828.412   # H122: 16668.375 -> 15828.412
$ws.Cells.Item(122, 9).Value2 = 18968.154   # I122: 20349.75 -> 18968.154
$ws.Cells.Item(122, 11).Value2 = 56904.462   # K122: 61049.25 -> 56904.462
$ws.Cells.Item(122, 13).Value2 = -54454.462   # M122: -58599.25 -> -54454.462
$ws.Cells.Item(126, 8).Value2 = 13255.12   # H126: 12445.148 -> 13255.12
$ws.Cells.Item(126, 9).Value2 = 18785.8   # I126: 20628.777 -> 18785.8
$ws.Cells.Item(126, 10).Value2 = 9568   # J126: 8353.333000000001 -> 9568
$ws.Cells.Item(126, 11).Value2 = 56357.39999999999   # K126: 61886.33099999999 -> 56357.39999999999
$ws.Cells.Item(126, 12).Value2 = 28704   # L126: 25059.999 -> 28704
$ws.Cells.Item(126, 13).Value2 = -53887.39999999999   # M126: -59416.33099999999 -> -53887.39999999999
$ws.Cells.Item(126, 14).Value2 = -33644   # N126: -29999.999 -> -33644
$ws.Cells.Item(134, 8).Value2 = 99900   # H134: 99450 -> 99900
$ws.Cells.Item(134, 10).Value2 = 99900   # J134: 99450 -> 99900
$ws.Cells.Item(134, 12).Value2 = 299700   # L134: 298350 -> 299700
$ws.Cells.Item(134, 14).Value2 = -304770   # N134: -303420 -> -304770
$ws.Cells.Item(135, 8).Value2 = 82182.836   # H135: 82308.39999999999 -> 82182.836
$ws.Cells.Item(135, 10).Value2 = 82182.836   # J135: 82308.39999999999 -> 82182.836
$ws.Cells.Item(135, 12).Value2 = 82182.836   # L135: 82308.39999999999 -> 82182.836
$ws.Cells.Item(135, 14).Value2 = -92322.836   # N135: -92448.39999999999 -> -92322.836

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 1897.1666   # H22: 1908.9445 -> 1897.1666
$ws.Cells.Item(22, 9).Value2 = 2003.2   # I22: 2193.7778 -> 2003.2
$ws.Cells.Item(22, 10).Value2 = 1764.625   # J22: 1624.1111 -> 1764.625
$ws.Cells.Item(22, 11).Value2 = 2003.2   # K22: 2193.7778 -> 2003.2
$ws.Cells.Item(22, 12).Value2 = 1764.625   # L22: 1624.1111 -> 1764.625
$ws.Cells.Item(22, 13).Value2 = -1708.2   # M22: -1898.7778 -> -1708.2
$ws.Cells.Item(22, 14).Value2 = -2354.625   # N22: -2214.1111 -> -2354.625
$ws.Cells.Item(27, 8).Value2 = 1897.1666   # H27: 1908.9445 -> 1897.1666
$ws.Cells.Item(27, 9).Value2 = 2003.2   # I27: 2193.7778 -> 2003.2
$ws.Cells.Item(27, 10).Value2 = 1764.625   # J27: 1624.1111 -> 1764.625
$ws.Cells.Item(27, 11).Value2 = 2003.2   # K27: 2193.7778 -> 2003.2
$ws.Cells.Item(27, 12).Value2 = 1764.625   # L27: 1624.1111 -> 1764.625
$ws.Cells.Item(27, 13).Value2 = -1896.2   # M27: -2086.7778 -> -1896.2
$ws.Cells.Item(27, 14).Value2 = -1978.625   # N27: -1838.1111 -> -1978.625
$ws.Cells.Item(40, 8).Value2 = 47488.35   # H40: 40954.95 -> 47488.35
$ws.Cells.Item(40, 9).Value2 = 58100.92   # I40: 47944.312 -> 58100.92
$ws.Cells.Item(40, 11).Value2 = 58100.92   # K40: 47944.312 -> 58100.92
$ws.Cells.Item(40, 13).Value2 = -57964.92   # M40: -47808.312 -> -57964.92
$ws.Cells.Item(55, 8).Value2 = 865.3200000000001   # H55: 822.1429000000001 -> 865.3200000000001
$ws.Cells.Item(55, 9).Value2 = 874.3   # I55: 840.0454999999999 -> 874.3
$ws.Cells.Item(55, 10).Value2 = 829.4   # J55: 756.5 -> 829.4
$ws.Cells.Item(55, 11).Value2 = 874.3   # K55: 840.0454999999999 -> 874.3
$ws.Cells.Item(55, 12).Value2 = 829.4   # L55: 756.5 -> 829.4
$ws.Cells.Item(55, 13).Value2 = -701.3   # M55: -667.0454999999999 -> -701.3
$ws.Cells.Item(55, 14).Value2 = -1175.4   # N55: -1102.5 -> -1175.4
$ws.Cells.Item(93, 8).Value2 = 4260.5557   # H93: 4083.6316 -> 4260.5557
$ws.Cells.Item(93, 9).Value2 = 4293.375   # I93: 4546.067 -> 4293.375
$ws.Cells.Item(93, 10).Value2 = 3998   # J93: 2349.5 -> 3998
$ws.Cells.Item(93, 11).Value2 = 4293.375   # K93: 4546.067 -> 4293.375
$ws.Cells.Item(93, 12).Value2 = 3998   # L93: 2349.5 -> 3998
$ws.Cells.Item(93, 13).Value2 = -3045.375   # M93: -3298.067 -> -3045.375
$ws.Cells.Item(93, 14).Value2 = -6494   # N93: -4845.5 -> -6494

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value2 = 23383.65   # H126: 19188.84 -> 23383.65
$ws.Cells.Item(126, 9).Value2 = 24315.723   # I126: 20939.572 -> 24315.723
$ws.Cells.Item(126, 10).Value2 = 14995   # J126: 9997.5 -> 14995
$ws.Cells.Item(126, 11).Value2 = 72947.16900000001   # K126: 62818.716 -> 72947.16900000001
$ws.Cells.Item(126, 12).Value2 = 44985   # L126: 29992.5 -> 44985
$ws.Cells.Item(126, 13).Value2 = -70477.16900000001   # M126: -60348.716 -> -70477.16900000001
$ws.Cells.Item(126, 14).Value2 = -49925   # N126: -34932.5 -> -49925
